# Organization and module adjustment
# A new "module" entry (key "unnamed" / Chinese value "待命名") is inserted
# into the module table on the "module" sheet, at row 273 (pushing the table
# rows below it down by one position). Two new entries are appended to the
# shared string table as a side effect of adding the new localized strings.
# A couple of rows further down the same column C list (rows 440/441, values
# "油井"/"煤层气井") also shift down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("module")

# --- Shift the ITEM/VALUE (columns B:C) block for rows 273-437 down to
# --- rows 274-438, making room for the newly-added row 273. Column A (the
# --- plain numeric index, already equal to row-1 for every populated row)
# --- is left untouched on purpose - it is not part of this shift.
$block = $ws.Range("B273:C437").Value()
$ws.Range("B274:C438").Value = $block

# New row 273: the freshly added localized entry.
$ws.Cells.Item(273, 2).Value = "unnamed"
$ws.Cells.Item(273, 3).Value = "待命名"

# Row 438 is now populated (it used to be a blank filler row) - give it the
# next sequential index to match the A = row-1 pattern used throughout the
# table.
$ws.Cells.Item(438, 1).Value = 437

# --- The same kind of single-row shift happens further down, in the
# --- isolated C440/C441 values ("油井" / "煤层气井"), which move to
# --- C441/C442.
$ws.Cells.Item(442, 3).Value = $ws.Cells.Item(441, 3).Value()
$ws.Cells.Item(441, 3).Value = $ws.Cells.Item(440, 3).Value()
$ws.Cells.Item(440, 3).Value = ""

# --- The sheet's used range grows by one row (A1:C642 -> A1:C643); extend
# --- it with a blank, styled row to match, by copying the last existing
# --- (blank) row down.
$ws.Range("A642:C642").Copy($ws.Range("A643:C643"))

# --- Reflect where the editor ended up after making this change (view
# --- scroll position + active selection).
$ws.Application.ActiveWindow.ScrollRow = 253
$ws.Range("B273").Select()
